$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text so Excel does not coerce '26.607.07' into a numeric value
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '26.607.07'
$ws.Range('E2').Value = '  -0.09%  '
# Force text so Excel does not coerce '1.595.38' into a numeric value
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.595.38'
$ws.Range('E3').Value = '  +0.48%  '
$ws.Range('E4').Value = '  +0.00%  '
# Force text so Excel does not coerce '211.46' into a numeric value
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '211.46'
$ws.Range('E5').Value = '  +0.38%  '
# Force text so Excel does not coerce '0.515' into a numeric value
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.515'
$ws.Range('E6').Value = '  +1.28%  '
$ws.Range('E7').Value = '  +0.01%  '
$ws.Range('E8').Value = '  +0.29%  '
$ws.Range('E9').Value = '  -0.54%  '
# Force text so Excel does not coerce '19.40' into a numeric value
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '19.40'
$ws.Range('E10').Value = '  -0.62%  '
$ws.Range('E11').Value = '  +0.47%  '
# Force text so Excel does not coerce '1.818.49' into a numeric value
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.818.49'
$ws.Range('E12').Value = '  +0.46%  '
# Force text so Excel does not coerce '1.619.14' into a numeric value
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '1.619.14'
$ws.Range('E13').Value = '  +2.03%  '
# Force text so Excel does not coerce '4.03' into a numeric value
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '4.03'
$ws.Range('E14').Value = '  +0.29%  '
$ws.Range('E15').Value = '  +0.17%  '
# Force text so Excel does not coerce '64.58' into a numeric value
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '64.58'
$ws.Range('E16').Value = '  -0.15%  '
# Force text so Excel does not coerce '26.590.66' into a numeric value
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '26.590.66'
$ws.Range('E17').Value = '  -0.09%  '
$ws.Range('E18').Value = '  +0.48%  '
# Force text so Excel does not coerce '208.67' into a numeric value
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '208.67'
$ws.Range('E19').Value = '  +0.41%  '
$ws.Range('E20').Value = '  +0.04%  '
# Force text so Excel does not coerce '6.96' into a numeric value
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '6.96'
$ws.Range('E21').Value = '  +3.86%  '
$ws.Range('E22').Value = '  +0.47%  '
$ws.Range('E23').Value = '  -1.84%  '
# Force text so Excel does not coerce '8.87' into a numeric value
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '8.87'
$ws.Range('E24').Value = '  +0.18%  '
# Force text so Excel does not coerce '145.31' into a numeric value
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '145.31'
$ws.Range('E25').Value = '  -0.86%  '
$ws.Range('E26').Value = '  -0.11%  '
$ws.Range('E27').Value = '  -1.59%  '
# Force text so Excel does not coerce '0.115' into a numeric value
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.115'
$ws.Range('E28').Value = '  +0.94%  '
# Force text so Excel does not coerce '15.25' into a numeric value
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '15.25'
$ws.Range('E29').Value = '  -0.04%  '
# Force text so Excel does not coerce '0.0506' into a numeric value
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.0506'
$ws.Range('E30').Value = '  -0.03%  '
$ws.Range('E31').Value = '  +0.45%  '
$ws.Range('E32').Value = '  +0.48%  '
# Force text so Excel does not coerce '0.653' into a numeric value
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.653'
$ws.Range('E33').Value = '  -0.40%  '
# Force text so Excel does not coerce '2.93' into a numeric value
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '2.93'
$ws.Range('E34').Value = '  +1.17%  '
# Force text so Excel does not coerce '1.281.15' into a numeric value
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.281.15'
$ws.Range('E35').Value = '  -1.73%  '
# Force text so Excel does not coerce '2.44' into a numeric value
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '2.44'
$ws.Range('E37').Value = '  +0.47%  '
$ws.Range('E38').Value = '  -0.03%  '
# Force text so Excel does not coerce '0.842' into a numeric value
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.842'
$ws.Range('E39').Value = '  +1.95%  '
$ws.Range('E40').Value = '  +0.02%  '
# Force text so Excel does not coerce '5.47' into a numeric value
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '5.47'
$ws.Range('E41').Value = '  +1.98%  '
$ws.Range('E42').Value = '  +1.19%  '
# Force text so Excel does not coerce '0.784' into a numeric value
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.784'
$ws.Range('E43').Value = '  -0.95%  '
# Force text so Excel does not coerce '64.03' into a numeric value
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '64.03'
$ws.Range('E44').Value = '  +2.46%  '
# Force text so Excel does not coerce '1.731.17' into a numeric value
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.731.17'
$ws.Range('E45').Value = '  +0.42%  '
# Force text so Excel does not coerce '0.912' into a numeric value
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.912'
$ws.Range('E46').Value = '  +8.87%  '
# Force text so Excel does not coerce '89.59' into a numeric value
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '89.59'
$ws.Range('E47').Value = '  +0.01%  '
$ws.Range('E48').Value = '  -1.02%  '
$ws.Range('E49').Value = '  -2.16%  '
$ws.Range('E50').Value = '  +4.10%  '
# Force text so Excel does not coerce '0.0506' into a numeric value
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.0506'
$ws.Range('E51').Value = '  +0.45%  '
